$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename sheet1 -> TC156
$ws.Name = "TC156"

# Fill in header + data (order matches shared-string pool order in target)
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"
$ws.Range("A1").Value = "Username"

# Apply thin box border around A1:B2
$rng = $ws.Range("A1:B2")
$rng.Borders.LineStyle = 1

# Set zoom and selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 240
$null = $ws.Range("A2").Select()
